$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FR")

$ws.Range("A29").Value = "Ecart moyen"
$ws.Range("B29").Formula = "=SUM(B2:B28)/22"

$ws.Range("A30").Value = "Somme des ecarts"
$ws.Range("B30").Formula = "=SUM(B2:B28)"
